$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values below are written as literal text (matching the source file's
# inlineStr cell type). Values that look numeric (e.g. "0.622", "10.91") are
# prefixed with a leading apostrophe so Excel stores them as text instead of
# re-interpreting/reformatting them as numbers.

$ws.Range("D2").Value = '34.099.28'
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").Value = '1.781.75'
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '''225.40'
$ws.Range("E5").Value = '  -0.77%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").Value = '''31.85'
$ws.Range("E8").Value = '  -1.27%  '
$ws.Range("E9").Value = '  -1.38%  '
$ws.Range("D10").Value = '''0.0687'
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("D11").Value = '''0.0948'
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("D12").Value = '2.038.31'
$ws.Range("E12").Value = '  -0.47%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = '''10.91'
$ws.Range("E13").Value = '  -4.72%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.770.18'
$ws.Range("E14").Value = '  -1.35%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '''0.622'
$ws.Range("E15").Value = '  -0.48%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '34.079.20'
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("D18").Value = '''67.57'
$ws.Range("E18").Value = '  -0.74%  '
$ws.Range("D19").Value = '''246.09'
$ws.Range("E19").Value = '  +1.31%  '
$ws.Range("D20").Value = '0.0₃0788'
$ws.Range("E20").Value = '  +1.65%  '
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("E22").Value = '  +0.87%  '
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("E24").Value = '  -0.81%  '
$ws.Range("D25").Value = '''161.91'
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("E26").Value = '  -0.76%  '
$ws.Range("D27").Value = '''16.27'
$ws.Range("E27").Value = '  +0.16%  '
$ws.Range("E28").Value = '  +0.39%  '
$ws.Range("E29").Value = '  +0.29%  '
$ws.Range("D30").Value = '''1.22'
$ws.Range("E30").Value = '  -1.44%  '
$ws.Range("E31").Value = '  -0.07%  '
$ws.Range("D32").Value = '''3.72'
$ws.Range("E32").Value = '  +1.62%  '
$ws.Range("E33").Value = '  +2.22%  '
$ws.Range("E34").Value = '  -2.54%  '
$ws.Range("D35").Value = '1.449.33'
$ws.Range("E35").Value = '  +2.92%  '
$ws.Range("E36").Value = '  +4.67%  '
$ws.Range("E37").Value = '  -0.51%  '
$ws.Range("E38").Value = '  +0.74%  '
$ws.Range("E39").Value = '  -0.88%  '
$ws.Range("E40").Value = '  +1.35%  '
$ws.Range("D41").Value = '''80.61'
$ws.Range("E41").Value = '  +0.71%  '
$ws.Range("E42").Value = '  +1.21%  '
$ws.Range("D43").Value = '''0.916'
$ws.Range("E43").Value = '  -0.79%  '
$ws.Range("E44").Value = '  +1.23%  '
$ws.Range("D45").Value = '''0.0519'
$ws.Range("E45").Value = '  +2.22%  '
$ws.Range("D46").Value = '''6.05'
$ws.Range("E46").Value = '  -0.75%  '
$ws.Range("E47").Value = '  +0.22%  '
$ws.Range("D48").Value = '1.938.32'
$ws.Range("E48").Value = '  -0.57%  '
$ws.Range("D49").Value = '0.0₆0132'
$ws.Range("E49").Value = '  -6.22%  '
$ws.Range("D50").Value = '''104.65'
$ws.Range("E50").Value = '  -2.86%  '
$ws.Range("E51").Value = '  +0.23%  '
